# Generacion de Ventana Multigraficos
# Adds the outgoing-mail configuration block (correo/pass/asunto/mensaje_1/mensaje_2/server)
# used to send "recuperar contraseña" emails, below the existing user-privilege table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: correo -------------------------------------------------------
$ws.Range("A9").Value = "correo"
# --- Row 10: pass ---------------------------------------------------------
$ws.Range("A10").Value = "pass"

# Column B for rows 9-10 (entered after the A column so the shared-string
# table order matches: correo, pass, email, app-password, ...)
$ws.Range("B9").Value = "eeganalysistoolbox@gmail.com"
$ws.Range("B10").Value = "cgsw pylb ptlf wvng"

# Turn the e-mail address into a real mailto: hyperlink (adds the
# built-in "Hyperlink" cell style + font).
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:eeganalysistoolbox@gmail.com")

# --- Row 11: asunto --------------------------------------------------------
$ws.Range("A11").Value = "asunto"
$ws.Range("B11").Value = "EEG Analysis Toolbox - Recuperar Contraseña"

# --- Row 12: mensaje_1 ------------------------------------------------------
$ws.Range("A12").Value = "mensaje_1"
# Leading apostrophe forces this text-looking-like-a-formula token to be
# stored as literal text with a quote prefix (quotePrefix="1").
$ws.Range("B12").Value = "'->"

# --- Row 13: mensaje_2 -------------------------------------------------------
$ws.Range("A13").Value = "mensaje_2"
$ws.Range("B13").Value = "Puedes cambiarla en cualquier momento en la pantalla de configuración despues de iniciar sesión."
$ws.Range("A13:B13").WrapText = $true
$ws.Range("A13:B13").RowHeight = 29.4

# --- Row 14: server -----------------------------------------------------------
$ws.Range("A14").Value = "server"
$ws.Range("B14").Value = "smtp.gmail.com"

# Restore the selection/scroll position to around the new block, like the
# author left the workbook after entering this data.
$ws.Range("A6").Select() | Out-Null
